$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the "Uncertainty (unc)" column (E) sample data: every affected
# row's corrupted numeric placeholder (leftover date-serial values like
# 46025/46024/1 from a broken numFmt) is replaced with the correct "Low"
# text rating, matching the other untouched rows (24, 28, 29, 31-35) that
# already hold Low/Moderate/High text.

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,25,26,27,30,36,37,38,39,40,41,42,43)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")
    $cell.Style = "Normal"
    $cell.Value = "Low"
}
